$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.751.39"
$ws.Range("E2").Value = "  -3.11%  "

$ws.Range("D3").Value = "3.474.33"
$ws.Range("E3").Value = "  -2.88%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.51%  "

$ws.Range("D7").Value = "3.472.61"
$ws.Range("E7").Value = "  -2.96%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -2.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.48%  "

$ws.Range("E11").Value = "  +2.53%  "

$ws.Range("E12").Value = "  -4.07%  "

$ws.Range("E13").Value = "  -4.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.35%  "

$ws.Range("D15").Value = "4.061.93"
$ws.Range("E15").Value = "  -2.87%  "

$ws.Range("D16").Value = "3.473.36"
$ws.Range("E16").Value = "  -2.83%  "

$ws.Range("D17").Value = "66.847.63"
$ws.Range("E17").Value = "  -3.25%  "

$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "439.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.62%  "

$ws.Range("E23").Value = "  -5.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.74%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "3.612.30"
$ws.Range("E26").Value = "  -2.83%  "

$ws.Range("E27").Value = "  -9.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.21%  "

$ws.Range("E29").Value = "  -7.98%  "

$ws.Range("E30").Value = "  -3.72%  "

$ws.Range("E31").Value = "  -6.69%  "

$ws.Range("E32").Value = "  -1.28%  "

$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.93%  "

$ws.Range("D36").Value = "3.464.22"
$ws.Range("E36").Value = "  -2.95%  "

$ws.Range("E37").Value = "  -7.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.30%  "

$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0891"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("E43").Value = "  -10.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "

$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.58%  "

$ws.Range("E48").Value = "  -9.22%  "

$ws.Range("E49").Value = "  -4.84%  "

$ws.Range("E50").Value = "  -8.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.986"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.92%  "
